$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A50").Value = "2025/26"
$ws.Range("B50").Value = "AFCON"
$ws.Range("C50").Value = "1z35p4iuhfxxdfaqjwzkqn2fo"

$ws.Range("C50").Select()
